# Update countries & provincias Spain
# Applies the data refresh captured in the commit "Update countries & provincias Spain".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp text (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 21 de Abril de 2020 a las 15:52"

# --- Row 4: Estados Unidos ---
$ws.Cells.Item(4, 2).Value = 792958   # B4 - Casos totales
$ws.Cells.Item(4, 3).Value = 199      # C4 - Nuevos casos
$ws.Cells.Item(4, 5).Value = 678038   # E4 - Recuperados
$ws.Cells.Item(4, 7).Value = 17       # G4 - Muertes hoy
$ws.Cells.Item(4, 8).Value = 42531    # H4 - Muertes

# --- Row 26: Arabia Saudita ---
$ws.Cells.Item(26, 6).Value = 81      # F26 - Casos criticos

# --- Row 41: Serbia ---
$ws.Cells.Item(41, 6).Value = 101     # F41 - Casos criticos

# --- Rows 72/73: Azerbaiyan and Nueva Zelanda swap order (Azerbaiyan now ranks above Nueva Zelanda) ---
$ws.Cells.Item(72, 1).Value = "Azerbaiyan"
$ws.Cells.Item(72, 2).Value = 1480
$ws.Cells.Item(72, 3).Value = 44
$ws.Cells.Item(72, 4).Value = 865
$ws.Cells.Item(72, 5).Value = 595
$ws.Cells.Item(72, 6).Value = 16
$ws.Cells.Item(72, 7).Value = 1
$ws.Cells.Item(72, 8).Value = 20

$ws.Cells.Item(73, 1).Value = "Nueva Zelanda"
$ws.Cells.Item(73, 2).Value = 1445
$ws.Cells.Item(73, 3).Value = 5
$ws.Cells.Item(73, 4).Value = 1006
$ws.Cells.Item(73, 5).Value = 426
$ws.Cells.Item(73, 6).Value = 3
$ws.Cells.Item(73, 7).Value = 1
$ws.Cells.Item(73, 8).Value = 13

# --- Row 117: Kenia ---
$ws.Cells.Item(117, 4).Value = 74     # D117 - Casos activos
$ws.Cells.Item(117, 5).Value = 208    # E117 - Recuperados

# --- Row 142: Liberia ---
$ws.Cells.Item(142, 2).Value = 101    # B142 - Casos totales
$ws.Cells.Item(142, 3).Value = 2      # C142 - Nuevos casos
$ws.Cells.Item(142, 5).Value = 86     # E142 - Recuperados

# --- Row 148: Maldivas ---
$ws.Cells.Item(148, 2).Value = 83     # B148 - Casos totales
$ws.Cells.Item(148, 3).Value = 14     # C148 - Nuevos casos
$ws.Cells.Item(148, 5).Value = 67     # E148 - Recuperados
$ws.Cells.Item(148, 6).Value = 2      # F148 - Casos criticos

# --- Rows 164/165: Sierra Leona and Macao swap order (Sierra Leona now ranks above Macao) ---
$ws.Cells.Item(164, 1).Value = "Sierra Leona"
$ws.Cells.Item(164, 2).Value = 50
$ws.Cells.Item(164, 3).Value = 7
$ws.Cells.Item(164, 4).Value = 6
$ws.Cells.Item(164, 5).Value = 44
$ws.Cells.Item(164, 6).Value = 0

$ws.Cells.Item(165, 1).Value = "Macao"
$ws.Cells.Item(165, 2).Value = 45
$ws.Cells.Item(165, 4).Value = 24
$ws.Cells.Item(165, 5).Value = 21
$ws.Cells.Item(165, 6).Value = 1
